try {
    [System.IO.File]::WriteAllText("C:\temp\test.txt", "hello")
    Write-Host "Write OK"
    $content = [System.IO.File]::ReadAllText("C:\temp\test.txt")
    Write-Host "Read: $content"
} catch {
    Write-Host "ERR: $_"
}
